$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 6624, 46070.95833333334),
    @(3, 6509, 46070.96875),
    @(4, 6461, 46070.97916666666),
    @(5, 6422, 46070.98958333334),
    @(6, 6467, 46071),
    @(7, 6388, 46071.01041666666),
    @(8, 6384, 46071.02083333334),
    @(9, 6324, 46071.03125),
    @(10, 6322, 46071.04166666666),
    @(11, 6258, 46071.05208333334),
    @(12, 6252, 46071.0625),
    @(13, 6246, 46071.07291666666),
    @(14, 6210, 46071.08333333334),
    @(15, 6194, 46071.09375),
    @(16, 6224, 46071.10416666666),
    @(17, 6243, 46071.11458333334),
    @(18, 6175, 46071.125),
    @(19, 6197, 46071.13541666666),
    @(20, 6206, 46071.14583333334),
    @(21, 6292, 46071.15625),
    @(22, 6327, 46071.16666666666),
    @(23, 6392, 46071.17708333334),
    @(24, 6434, 46071.1875),
    @(25, 6523, 46071.19791666666),
    @(26, 6664, 46071.20833333334),
    @(27, 6850, 46071.21875),
    @(28, 6964, 46071.22916666666),
    @(29, 7102, 46071.23958333334),
    @(30, 7322, 46071.25),
    @(31, 7478, 46071.26041666666),
    @(32, 7572, 46071.27083333334),
    @(33, 7739, 46071.28125),
)

foreach ($row in $data) {
    $r = $row[0]
    $a = $row[1]
    $b = $row[2]
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
}

